$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "TEST_TRACKING_DECEMBER"
$ws.Range("B3").Value = "dhl"
$ws.Range("C3").Value = "Delivered"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "2025-12-31T23:59:59"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "2026-02-08T11:24:55+00:00"
$ws.Range("I3").Value = "TEST_Tracking_December"

# Row 4
$ws.Range("A4").Value = "TEST_TRACKING"
$ws.Range("B4").Value = "kn"
$ws.Range("C4").Value = "Delivered"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "2026-02-07T23:59:59"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "2026-02-08T11:24:18+00:00"
$ws.Range("I4").Value = "ZFRE"

# Row 5
$ws.Range("A5").Value = "ITD-0-12345678"
$ws.Range("B5").Value = "testing-courier"
$ws.Range("C5").Value = "Delivered"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "2026-02-08T07:07:24-04:00"
$ws.Range("G5").Value = "1000 W Test Ln, Austin, Texas, 11111"
$ws.Range("H5").Value = "2026-02-08T11:07:24+00:00"
$ws.Range("I5").Value = "ITD-0-12345678"
